$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "havainnointijaksot vuonna Orionin tähtikuvio 2022",
    $false,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "Orionin tähtikuvio havainnointijaksot vuonna 2022",
    2
)
